# DbSchema NSKFDC-App.xlsx : add the "GenerateReports" table as a new
# column (column I) on Sheet1, mirroring the layout of the other schema
# tables already on the sheet (User / Candidate / BankDetails / ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: clone the format of a reference cell onto a target cell,
# then write the target's value -------------------------------------------
function Set-Cell($row, $col, $styleRow, $styleCol, $value) {
    $ws.Cells.Item($styleRow, $styleCol).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, $col).Value = $value
}

# Table header "GenerateReports" -> same look as the other table headers
# (A2 / E2 / G2), yellow-filled bold box.
Set-Cell 2 9 2 7 "GenerateReports"

# Primary key row -> same look as the other "- PK" rows (A3), blue text
# with a border, plus a light header fill for the PK row.
Set-Cell 3 9 3 1 "generateReportsId (Auto Increment )- PK"
$ws.Cells.Item(3, 9).Interior.ThemeColor = 2

# Ordinary attribute rows -> same look as the other plain rows (A4).
Set-Cell 4 9 4 1 "occupationReport - Flag"
Set-Cell 5 9 4 1 "ORGeneratedOn"
Set-Cell 6 9 4 1 "attendanceSheet - Flag"
Set-Cell 7 9 4 1 "ASGeneratedOn"
Set-Cell 8 9 4 1 "NSKFDCSheet - Flag"
Set-Cell 9 9 4 1 "NSKFDCGeneratedOn"
Set-Cell 10 9 4 1 "SDMSSheet - Flag"
Set-Cell 11 9 4 1 "SDMSGeneratedOn"
Set-Cell 12 9 4 1 "selectionCommittee - Flag "
Set-Cell 13 9 4 1 "SCGeneratedOn"

# Last two attribute rows -> distinct (bordered, no-fill) look.
Set-Cell 14 9 4 1 "batchReport - Flag"
Set-Cell 15 9 4 1 "BRGeneratedOn"

# Trailing FK rows -> same look as the other "- FK" rows (C29).
Set-Cell 16 9 29 3 "batchId - FK"
Set-Cell 17 9 29 3 "trainingPartnerEmail - FK"

# A couple of stray single-space cells in column H that came along with
# the new column in the original edit.
$ws.Cells.Item(7, 8).Value = " "
$ws.Cells.Item(11, 8).Value = " "

# Widen column I to fit the new, longer table entries and scroll the
# sheet so the new column is in view, matching the saved view state.
$ws.Columns.Item(9).ColumnWidth = 37.42578125
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("I16").Select()
